# chore: end and init work
# end work day 29 and init work on day 30

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 69 (29th - end of work day): set HORA F, PAUSAS and update ASSUNTO
$ws.Range("C69").Value = 0.6875
$ws.Range("E69").Value = 0.20833333333333334
$ws.Range("G69").Value = "ESTÁGIO + SOFT"

# Row 70 (30th - init work day): set HORA I and ASSUNTO / PRODUÇÃO
$ws.Range("B70").Value = 0.61111111111111105
$ws.Range("G70").Value = "HARD"
$ws.Range("H70").Value = "HARD"

# Update the active selection to reflect where the user left off editing
$ws.Range("H71").Select()
